$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.3934994375533165 ; $ws.Range("C2").Value = 0.05307076781402031 ; $ws.Range("E2").Value = 0.6441152678879121 ; $ws.Range("F2").Value = 2.35707244740486 ; $ws.Range("G2").Value = 0.002446762567082351 ; $ws.Range("J2").Value = 0.0532035366023873 ; $ws.Range("K2").Value = 0.3699106016603082 ; $ws.Range("O2").Value = 2.727889815261506
$ws.Range("B3").Value = 0.351916771976164 ; $ws.Range("C3").Value = 0.04771218570084557 ; $ws.Range("E3").Value = 0.6139038305324931 ; $ws.Range("F3").Value = 2.331975730226617 ; $ws.Range("G3").Value = 0.002449180088866428 ; $ws.Range("J3").Value = 0.05391904622141475 ; $ws.Range("K3").Value = 0.3256463336316529 ; $ws.Range("O3").Value = 2.763607967927598
$ws.Range("B4").Value = 0.3263971591446762 ; $ws.Range("C4").Value = 0.04440530873365844 ; $ws.Range("E4").Value = 0.5956365294616575 ; $ws.Range("F4").Value = 2.318030432657693 ; $ws.Range("G4").Value = 0.002450741523717852 ; $ws.Range("J4").Value = 0.0543949113762352 ; $ws.Range("K4").Value = 0.2984228377471254 ; $ws.Range("O4").Value = 2.787404497879209
$ws.Range("B5").Value = 0.3160013584277408 ; $ws.Range("C5").Value = 0.04305361655693218 ; $ws.Range("E5").Value = 0.588263496618282 ; $ws.Range("F5").Value = 2.312715200856942 ; $ws.Range("G5").Value = 0.002451397258632465 ; $ws.Range("J5").Value = 0.05459800313517071 ; $ws.Range("K5").Value = 0.2873182993557748 ; $ws.Range("O5").Value = 2.797570180278356
$ws.Range("B6").Value = 0.314275380375733 ; $ws.Range("C6").Value = 0.04282892314121511 ; $ws.Range("E6").Value = 0.5870434999833947 ; $ws.Range("F6").Value = 2.311854795039494 ; $ws.Range("G6").Value = 0.002451507318692592 ; $ws.Range("J6").Value = 0.05463227993802455 ; $ws.Range("K6").Value = 0.2854737677618289 ; $ws.Range("O6").Value = 2.799286454687305
$ws.Range("B7").Value = 0.3262569422313106 ; $ws.Range("C7").Value = 0.04438709588403356 ; $ws.Range("E7").Value = 0.5955368065616824 ; $ws.Range("F7").Value = 2.317957262013309 ; $ws.Range("G7").Value = 0.002450750288191043 ; $ws.Range("J7").Value = 0.05439761322216796 ; $ws.Range("K7").Value = 0.2982731207742688 ; $ws.Range("O7").Value = 2.787539700186031
$ws.Range("B8").Value = 0.3791594966565413 ; $ws.Range("C8").Value = 0.0512266397377914 ; $ws.Range("E8").Value = 0.6336396706324763 ; $ws.Range("F8").Value = 2.348114835329127 ; $ws.Range("G8").Value = 0.002447580168845836 ; $ws.Range("J8").Value = 0.05344265141940951 ; $ws.Range("K8").Value = 0.3546579708228705 ; $ws.Range("O8").Value = 2.739817856981162
$ws.Range("B9").Value = 0.4829799351592214 ; $ws.Range("C9").Value = 0.06450357847303678 ; $ws.Range("E9").Value = 0.71060800579491 ; $ws.Range("F9").Value = 2.418904662693976 ; $ws.Range("G9").Value = 0.002441972376475143 ; $ws.Range("J9").Value = 0.05186057849647341 ; $ws.Range("K9").Value = 0.4648507062620126 ; $ws.Range("O9").Value = 2.66106708303181
$ws.Range("B10").Value = 0.559287041180113 ; $ws.Range("C10").Value = 0.0741725456347524 ; $ws.Range("E10").Value = 0.7685426300641751 ; $ws.Range("F10").Value = 2.478072130460845 ; $ws.Range("G10").Value = 0.00243821969425942 ; $ws.Range("J10").Value = 0.05087632085884053 ; $ws.Range("K10").Value = 0.5455598615515953 ; $ws.Range("O10").Value = 2.612291948717257
$ws.Range("B11").Value = 0.5940044568458802 ; $ws.Range("C11").Value = 0.07855204227205093 ; $ws.Range("E11").Value = 0.7952032494922179 ; $ws.Range("F11").Value = 2.506555886483 ; $ws.Range("G11").Value = 0.002436591465812832 ; $ws.Range("J11").Value = 0.0504674292287195 ; $ws.Range("K11").Value = 0.5822188870467642 ; $ws.Range("O11").Value = 2.592084108353603
$ws.Range("B12").Value = 0.6071513085961158 ; $ws.Range("C12").Value = 0.08020764812521008 ; $ws.Range("E12").Value = 0.8053430775790389 ; $ws.Range("F12").Value = 2.517568287345853 ; $ws.Range("G12").Value = 0.002435986180899631 ; $ws.Range("J12").Value = 0.05031819685039629 ; $ws.Range("K12").Value = 0.5960921893260434 ; $ws.Range("O12").Value = 2.584717443321807
$ws.Range("B13").Value = 0.6043199028459867 ; $ws.Range("C13").Value = 0.07985121033846099 ; $ws.Range("E13").Value = 0.8031573251430331 ; $ws.Range("F13").Value = 2.515186497055453 ; $ws.Range("G13").Value = 0.002436116038555208 ; $ws.Range("J13").Value = 0.05035008707914201 ; $ws.Range("K13").Value = 0.5931047211965961 ; $ws.Range("O13").Value = 2.586291270456883
$ws.Range("B14").Value = 0.595086057207908 ; $ws.Range("C14").Value = 0.07868830686463468 ; $ws.Range("E14").Value = 0.7960365765422921 ; $ws.Range("F14").Value = 2.507457345395778 ; $ws.Range("G14").Value = 0.002436541442681481 ; $ws.Range("J14").Value = 0.05045503932826989 ; $ws.Range("K14").Value = 0.5833604294042232 ; $ws.Range("O14").Value = 2.591472319535413
$ws.Range("B15").Value = 0.5894300583365748 ; $ws.Range("C15").Value = 0.07797562580515205 ; $ws.Range("E15").Value = 0.7916806478840925 ; $ws.Range("F15").Value = 2.502752498854136 ; $ws.Range("G15").Value = 0.00243680348360206 ; $ws.Range("J15").Value = 0.05052005624829725 ; $ws.Range("K15").Value = 0.5773906210746986 ; $ws.Range("O15").Value = 2.594683081314429
$ws.Range("B16").Value = 0.557018225177984 ; $ws.Range("C16").Value = 0.07388594612054078 ; $ws.Range("E16").Value = 0.7668064546410847 ; $ws.Range("F16").Value = 2.476242267065089 ; $ws.Range("G16").Value = 0.002438327685676725 ; $ws.Range("J16").Value = 0.05090382632516643 ; $ws.Range("K16").Value = 0.5431629280486447 ; $ws.Range("O16").Value = 2.61365249847448
$ws.Range("B17").Value = 0.537135457671809 ; $ws.Range("C17").Value = 0.07137214028475114 ; $ws.Range("E17").Value = 0.7516253150162555 ; $ws.Range("F17").Value = 2.460381256083394 ; $ws.Range("G17").Value = 0.002439282900346214 ; $ws.Range("J17").Value = 0.05114922137603095 ; $ws.Range("K17").Value = 0.5221505636697827 ; $ws.Range("O17").Value = 2.62579740832247
$ws.Range("B18").Value = 0.5256999146991745 ; $ws.Range("C18").Value = 0.06992448574285959 ; $ws.Range("E18").Value = 0.7429223117948851 ; $ws.Range("F18").Value = 2.451405983106397 ; $ws.Range("G18").Value = 0.002439839742839227 ; $ws.Range("J18").Value = 0.05129402191513321 ; $ws.Range("K18").Value = 0.5100595759998043 ; $ws.Range("O18").Value = 2.632969175136211
$ws.Range("B19").Value = 0.5218281400639455 ; $ws.Range("C19").Value = 0.06943403176295249 ; $ws.Range("E19").Value = 0.7399805672336583 ; $ws.Range("F19").Value = 2.448392431850195 ; $ws.Range("G19").Value = 0.002440029557582818 ; $ws.Range("J19").Value = 0.05134367619492508 ; $ws.Range("K19").Value = 0.5059648959805827 ; $ws.Range("O19").Value = 2.635429389551462
$ws.Range("B20").Value = 0.5392519648487166 ; $ws.Range("C20").Value = 0.07163992405710928 ; $ws.Range("E20").Value = 0.7532383935985081 ; $ws.Range("F20").Value = 2.462054410572136 ; $ws.Range("G20").Value = 0.002439180447715045 ; $ws.Range("J20").Value = 0.05112272015454344 ; $ws.Range("K20").Value = 0.5243879117144274 ; $ws.Range("O20").Value = 2.624485270683294
$ws.Range("B21").Value = 0.5977982629240159 ; $ws.Range("C21").Value = 0.07902995688002079 ; $ws.Range("E21").Value = 0.7981269176615058 ; $ws.Range("F21").Value = 2.509721439277911 ; $ws.Range("G21").Value = 0.002436416185264485 ; $ws.Range("J21").Value = 0.05042406001088828 ; $ws.Range("K21").Value = 0.586222804512289 ; $ws.Range("O21").Value = 2.589942761770686
$ws.Range("B22").Value = 0.6360620816365099 ; $ws.Range("C22").Value = 0.08384332570690844 ; $ws.Range("E22").Value = 0.8277208285300901 ; $ws.Range("F22").Value = 2.542193522584625 ; $ws.Range("G22").Value = 0.002434675362012983 ; $ws.Range("J22").Value = 0.05000013025101424 ; $ws.Range("K22").Value = 0.6265846295330562 ; $ws.Range("O22").Value = 2.56903235995334
$ws.Range("B23").Value = 0.6156401158140454 ; $ws.Range("C23").Value = 0.08127587350520571 ; $ws.Range("E23").Value = 0.8119025000781619 ; $ws.Range("F23").Value = 2.524741640042805 ; $ws.Range("G23").Value = 0.002435598470629796 ; $ws.Range("J23").Value = 0.050223392521346 ; $ws.Range("K23").Value = 0.6050476224430668 ; $ws.Range("O23").Value = 2.580039993668649
$ws.Range("B24").Value = 0.5382951064727592 ; $ws.Range("C24").Value = 0.07151886658874673 ; $ws.Range("E24").Value = 0.7525090434157846 ; $ws.Range("F24").Value = 2.461297530686068 ; $ws.Range("G24").Value = 0.00243922674266103 ; $ws.Range("J24").Value = 0.05113468977202729 ; $ws.Range("K24").Value = 0.523376439945423 ; $ws.Range("O24").Value = 2.625077897890051
$ws.Range("B25").Value = 0.4548871821866385 ; $ws.Range("C25").Value = 0.06092662898531387 ; $ws.Range("E25").Value = 0.6895437410525318 ; $ws.Range("F25").Value = 2.398501114096916 ; $ws.Range("G25").Value = 0.002443424648553187 ; $ws.Range("J25").Value = 0.05225736398586811 ; $ws.Range("K25").Value = 0.4350830942998414 ; $ws.Range("O25").Value = 2.680779163078569
